$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($text -eq "redux") {
        # wdGreen = 4 : highlight both the run and the paragraph mark
        # so the list item ("redux") is fully green-highlighted.
        $p.Range.Font.HighlightColorIndex = 4
    }
}
